# SB response further process
# Apply changes to the "Amazon_HIT_expert_submission_track_master" sheet:
#  - Row 7 (the DW "pairwise_merge" summary row) is removed; its merged-file
#    note is moved into the new "Further_process" (F) column of row 6.
#  - Two new "Further_process" notes are recorded for the WS/Wali and the
#    SB rows.
#  - The old last row (SB "video does not play" resub tracker row) is removed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move the merge-summary note from row 7 up into the new F column on row 6,
# then delete row 7 entirely (its A:E content is now redundant).
$ws.Range("F6").Value = $ws.Range("E7").Value2
$ws.Rows(7).Delete()

# Record further-processing notes for the Wali (Jun-19 to Jul-14) merge.
$ws.Range("F9").Value = "master_all_responses_Jun-19-2023_to_Jul-14-2023_Wali.csv"

# Record further-processing notes for the SB (Oct-01 to SB resub Oct-01) merge.
$ws.Range("F15").Value = "master_all_responses_SB_Oct-01-2023_to_SB_resub_Oct-01-2023_Sarah.csv"

# Remove the trailing SB "video does not play" resub-tracker row (old row 17,
# now row 16 after the earlier deletion above).
$ws.Rows(16).Delete()

$ws.Range("D14").Select()
